$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Price column (D) updates: force text format to preserve exact
#     string representation (trailing zeros, etc.), matching the
#     original inline-string cell type.
$priceCells = @("D2", "D3", "D4", "D5", "D6", "D7", "D8", "D9", "D10", "D11", "D12", "D14", "D15", "D16", "D17", "D18", "D19", "D20", "D21", "D22", "D23", "D24", "D25", "D26", "D40", "D41", "D42", "D43", "D44", "D45", "D48")
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "274.47"
$ws.Range("D3").Value = "22.99"
$ws.Range("D4").Value = "6.362"
$ws.Range("D5").Value = "0.06236"
$ws.Range("D6").Value = "3.658"
$ws.Range("D7").Value = "6.685"
$ws.Range("D8").Value = "1.368"
$ws.Range("D9").Value = "0.8313"
$ws.Range("D10").Value = "0.01380"
$ws.Range("D11").Value = "0.1637"
$ws.Range("D12").Value = "0.08273"
$ws.Range("D14").Value = "0.03105"
$ws.Range("D15").Value = "0.09312"
$ws.Range("D16").Value = "3.882"
$ws.Range("D17").Value = "0.001669"
$ws.Range("D18").Value = "0.04772"
$ws.Range("D19").Value = "0.006340"
$ws.Range("D20").Value = "0.005562"
$ws.Range("D21").Value = "0.001089"
$ws.Range("D22").Value = "0.0001499"
$ws.Range("D23").Value = "3.725"
$ws.Range("D24").Value = "2.370"
$ws.Range("D25").Value = "0.3380"
$ws.Range("D26").Value = "0.1270"
$ws.Range("D40").Value = "0.04692"
$ws.Range("D41").Value = "0.007041"
$ws.Range("D42").Value = "0.1165"
$ws.Range("D43").Value = "0.003599"
$ws.Range("D44").Value = "0.01217"
$ws.Range("D45").Value = "0.00006255"
$ws.Range("D48").Value = "0.03206"

# --- Text column (B, C, E) updates: coin names, links, labels.
$ws.Range("B15").Value = "BitMartToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("E15").Value = "14BitMartTokenBMX"
$ws.Range("B16").Value = "MCDex"
$ws.Range("C16").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("E16").Value = "15MCDexMCB"
$ws.Range("B17").Value = "BitForexToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("E17").Value = "16BitForexTokenBF"
$ws.Range("B18").Value = "CoinExToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("E18").Value = "17CoinExTokenCET"
$ws.Range("B19").Value = "TigerCash"
$ws.Range("C19").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("E19").Value = "18TigerCashTCH"
$ws.Range("B20").Value = "HotbitToken"
$ws.Range("C20").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
$ws.Range("E20").Value = "19HotbitTokenHTBWorstin24h"
$ws.Range("B21").Value = "BitKan"
$ws.Range("C21").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
$ws.Range("E21").Value = "20BitKanKAN"
$ws.Range("B22").Value = "NitroEx"
$ws.Range("C22").Value = "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"
$ws.Range("E22").Value = "21NitroExNTX"
$ws.Range("B23").Value = "LEO"
$ws.Range("C23").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("E23").Value = "22LEOLEO"
$ws.Range("B24").Value = "BTSEToken"
$ws.Range("C24").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("E24").Value = "23BTSETokenBTSE"
$ws.Range("B25").Value = "BitpandaEcosystemToken"
$ws.Range("C25").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("E25").Value = "24BitpandaEcosystemTokenBEST"
$ws.Range("B26").Value = "ProBitToken"
$ws.Range("C26").Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
$ws.Range("E26").Value = "25ProBitTokenPROB"
$ws.Range("E43").Value = "42CEJICEJI"
